$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# The paragraph that currently holds the bold "Play Break Away Deluxe..."
# text (near the end of the doc, second-to-last paragraph) already has the
# exact run layout we want to reproduce here (a leading empty run followed
# by a bold run) -- copy its formatted content over, then swap the wording
# in place. This keeps the resulting markup consistent with how Word
# itself represents "leading empty run + bold run" paragraphs in this doc.
$titleBoldPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$metaPara.Range.FormattedText = $titleBoldPara.Range.FormattedText

$metaFind = $metaPara.Range.Find
$metaFind.Execute("Play Break Away Deluxe Free - Unique 5x5 Grid Slot Game", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Meta description", 2) | Out-Null

# Append the (non-bold) remainder of the meta description right after the
# bold "Meta description" label.
$afterLabel = $d.Range($metaPara.Range.End - 1, $metaPara.Range.End - 1)
$afterLabel.InsertAfter(": Read our review of Break Away Deluxe, a unique Microgaming slot game with a 5x5 grid structure, falling symbols function, and Smashing Wild feature. Play for free now.")
$afterLabel.Bold = 0

# ---------------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Break Away Deluxe..." paragraph that
#    used to sit just before the closing italic paragraph. (Its index moved
#    from 49 -> 50 because of the paragraph inserted near the top above.)
# ---------------------------------------------------------------------------

$dupTitlePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$dupTitlePara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the closing italic paragraph with the new image
#    prompt, keeping its italic run formatting intact.
# ---------------------------------------------------------------------------

$closingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$closingTextRange = $d.Range($closingPara.Range.Start, $closingPara.Range.End - 1)
$closingTextRange.Text = "Create a feature image for ""Break Away Deluxe"" that features a happy Maya warrior with glasses in a cartoon style. The image should showcase the excitement of playing the game and capture the theme of ice hockey in a fun and playful way. The warrior should be wearing ice skates and holding a hockey stick, with the game's logo in the background. The colors used should be bright and bold to grab the viewer's attention and make them want to try the game. The overall design should convey the thrill of playing and winning on the ice, while also incorporating the unique element of the Maya warrior for added interest."

Write-Output "done"
